$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 16:22"

# Update country rows whose data shifted (country name + stats)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 850078
$ws.Cells.Item(4, 3).Value = 1361
$ws.Cells.Item(4, 4).Value = 84058
$ws.Cells.Item(4, 5).Value = 718286
$ws.Cells.Item(4, 6).Value = 14016
$ws.Cells.Item(4, 7).Value = 75
$ws.Cells.Item(4, 8).Value = 47734

$ws.Cells.Item(8, 1).Value = "Alemania"
$ws.Cells.Item(8, 2).Value = 151022
$ws.Cells.Item(8, 3).Value = 374
$ws.Cells.Item(8, 4).Value = 103300
$ws.Cells.Item(8, 5).Value = 42368
$ws.Cells.Item(8, 6).Value = 2908
$ws.Cells.Item(8, 7).Value = 39
$ws.Cells.Item(8, 8).Value = 5354

$ws.Cells.Item(9, 1).Value = "Reino Unido"
$ws.Cells.Item(9, 2).Value = 138078
$ws.Cells.Item(9, 3).Value = 4583
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 118996
$ws.Cells.Item(9, 6).Value = 1559
$ws.Cells.Item(9, 7).Value = 638
$ws.Cells.Item(9, 8).Value = 18738

$ws.Cells.Item(17, 1).Value = "Paises Bajos"
$ws.Cells.Item(17, 2).Value = 35729
$ws.Cells.Item(17, 3).Value = 887
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 31302
$ws.Cells.Item(17, 6).Value = 1008
$ws.Cells.Item(17, 7).Value = 123
$ws.Cells.Item(17, 8).Value = 4177

$ws.Cells.Item(22, 1).Value = "Suecia"
$ws.Cells.Item(22, 2).Value = 16755
$ws.Cells.Item(22, 3).Value = 751
$ws.Cells.Item(22, 4).Value = 550
$ws.Cells.Item(22, 5).Value = 14184
$ws.Cells.Item(22, 6).Value = 533
$ws.Cells.Item(22, 7).Value = 84
$ws.Cells.Item(22, 8).Value = 2021

$ws.Cells.Item(26, 1).Value = "Arabia Saudita"
$ws.Cells.Item(26, 2).Value = 13930
$ws.Cells.Item(26, 3).Value = 1158
$ws.Cells.Item(26, 4).Value = 1925
$ws.Cells.Item(26, 5).Value = 11884
$ws.Cells.Item(26, 6).Value = 82
$ws.Cells.Item(26, 7).Value = 7
$ws.Cells.Item(26, 8).Value = 121

$ws.Cells.Item(35, 1).Value = "Rumania"
$ws.Cells.Item(35, 2).Value = 10096
$ws.Cells.Item(35, 3).Value = 386
$ws.Cells.Item(35, 4).Value = 2478
$ws.Cells.Item(35, 5).Value = 7077
$ws.Cells.Item(35, 6).Value = 236
$ws.Cells.Item(35, 7).Value = 17
$ws.Cells.Item(35, 8).Value = 541

$ws.Cells.Item(38, 1).Value = "Bielorrusia"
$ws.Cells.Item(38, 2).Value = 8022
$ws.Cells.Item(38, 3).Value = 741
$ws.Cells.Item(38, 4).Value = 938
$ws.Cells.Item(38, 5).Value = 7024
$ws.Cells.Item(38, 6).Value = 92
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = 60

$ws.Cells.Item(39, 1).Value = "Indonesia"
$ws.Cells.Item(39, 2).Value = 7775
$ws.Cells.Item(39, 3).Value = 357
$ws.Cells.Item(39, 4).Value = 960
$ws.Cells.Item(39, 5).Value = 6168
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 12
$ws.Cells.Item(39, 8).Value = 647

$ws.Cells.Item(40, 1).Value = "Catar"
$ws.Cells.Item(40, 2).Value = 7764
$ws.Cells.Item(40, 3).Value = 623
$ws.Cells.Item(40, 4).Value = 750
$ws.Cells.Item(40, 5).Value = 7004
$ws.Cells.Item(40, 6).Value = 72
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 10

$ws.Cells.Item(41, 1).Value = "Noruega"
$ws.Cells.Item(41, 2).Value = 7361
$ws.Cells.Item(41, 3).Value = 23
$ws.Cells.Item(41, 4).Value = 32
$ws.Cells.Item(41, 5).Value = 7138
$ws.Cells.Item(41, 6).Value = 54
$ws.Cells.Item(41, 7).Value = 4
$ws.Cells.Item(41, 8).Value = 191

$ws.Cells.Item(42, 1).Value = "Serbia"
$ws.Cells.Item(42, 2).Value = 7276
$ws.Cells.Item(42, 3).Value = 162
$ws.Cells.Item(42, 4).Value = 1063
$ws.Cells.Item(42, 5).Value = 6074
$ws.Cells.Item(42, 6).Value = 96
$ws.Cells.Item(42, 7).Value = 5
$ws.Cells.Item(42, 8).Value = 139

$ws.Cells.Item(43, 1).Value = "Ucrania"
$ws.Cells.Item(43, 2).Value = 7170
$ws.Cells.Item(43, 3).Value = 578
$ws.Cells.Item(43, 4).Value = 504
$ws.Cells.Item(43, 5).Value = 6479
$ws.Cells.Item(43, 6).Value = 45
$ws.Cells.Item(43, 7).Value = 13
$ws.Cells.Item(43, 8).Value = 187

$ws.Cells.Item(44, 1).Value = "Chequia"
$ws.Cells.Item(44, 2).Value = 7136
$ws.Cells.Item(44, 3).Value = 4
$ws.Cells.Item(44, 4).Value = 2002
$ws.Cells.Item(44, 5).Value = 4924
$ws.Cells.Item(44, 6).Value = 76
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(44, 8).Value = 210

$ws.Cells.Item(57, 1).Value = "Argentina"
$ws.Cells.Item(57, 2).Value = 3288
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(57, 4).Value = 919
$ws.Cells.Item(57, 5).Value = 2210
$ws.Cells.Item(57, 6).Value = 123
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 159

$ws.Cells.Item(58, 1).Value = "Moldavia"
$ws.Cells.Item(58, 2).Value = 2926
$ws.Cells.Item(58, 3).Value = 148
$ws.Cells.Item(58, 4).Value = 661
$ws.Cells.Item(58, 5).Value = 2186
$ws.Cells.Item(58, 6).Value = 212
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 79

$ws.Cells.Item(59, 1).Value = "Argelia"
$ws.Cells.Item(59, 2).Value = 2910
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 1204
$ws.Cells.Item(59, 5).Value = 1304
$ws.Cells.Item(59, 6).Value = 40
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 402

$ws.Cells.Item(60, 1).Value = "Tailandia"
$ws.Cells.Item(60, 2).Value = 2839
$ws.Cells.Item(60, 3).Value = 13
$ws.Cells.Item(60, 4).Value = 2430
$ws.Cells.Item(60, 5).Value = 359
$ws.Cells.Item(60, 6).Value = 61
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 50

$ws.Cells.Item(64, 1).Value = "Kazajistan"
$ws.Cells.Item(64, 2).Value = 2251
$ws.Cells.Item(64, 3).Value = 116
$ws.Cells.Item(64, 4).Value = 551
$ws.Cells.Item(64, 5).Value = 1680
$ws.Cells.Item(64, 6).Value = 29
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = 20

$ws.Cells.Item(67, 1).Value = "Islandia"
$ws.Cells.Item(67, 2).Value = 1789
$ws.Cells.Item(67, 3).Value = 4
$ws.Cells.Item(67, 4).Value = 1509
$ws.Cells.Item(67, 5).Value = 270
$ws.Cells.Item(67, 6).Value = 5
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 10

$ws.Cells.Item(69, 1).Value = "Uzbekistan"
$ws.Cells.Item(69, 2).Value = 1716
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 4).Value = 495
$ws.Cells.Item(69, 5).Value = 1214
$ws.Cells.Item(69, 6).Value = 8
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 7

$ws.Cells.Item(80, 1).Value = "Afganistan"
$ws.Cells.Item(80, 2).Value = 1279
$ws.Cells.Item(80, 3).Value = 103
$ws.Cells.Item(80, 4).Value = 179
$ws.Cells.Item(80, 5).Value = 1058
$ws.Cells.Item(80, 6).Value = 7
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = 42

$ws.Cells.Item(86, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(86, 2).Value = 986
$ws.Cells.Item(86, 3).Value = 12
$ws.Cells.Item(86, 4).Value = 252
$ws.Cells.Item(86, 5).Value = 732
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 2

$ws.Cells.Item(104, 1).Value = "San Marino"
$ws.Cells.Item(104, 2).Value = 501
$ws.Cells.Item(104, 3).Value = 13
$ws.Cells.Item(104, 4).Value = 63
$ws.Cells.Item(104, 5).Value = 398
$ws.Cells.Item(104, 6).Value = 3
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 40

$ws.Cells.Item(114, 1).Value = "Sri Lanka"
$ws.Cells.Item(114, 2).Value = 335
$ws.Cells.Item(114, 3).Value = 5
$ws.Cells.Item(114, 4).Value = 107
$ws.Cells.Item(114, 5).Value = 221
$ws.Cells.Item(114, 6).Value = 2
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 7

$ws.Cells.Item(117, 1).Value = "Kenia"
$ws.Cells.Item(117, 2).Value = 320
$ws.Cells.Item(117, 3).Value = 17
$ws.Cells.Item(117, 4).Value = 89
$ws.Cells.Item(117, 5).Value = 217
$ws.Cells.Item(117, 6).Value = 2
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 14

$ws.Cells.Item(118, 1).Value = "Montenegro"
$ws.Cells.Item(118, 2).Value = 316
$ws.Cells.Item(118, 3).Value = 1
$ws.Cells.Item(118, 4).Value = 116
$ws.Cells.Item(118, 5).Value = 195
$ws.Cells.Item(118, 6).Value = 7
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 5

$ws.Cells.Item(119, 1).Value = "Isla de Man"
$ws.Cells.Item(119, 2).Value = 307
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 212
$ws.Cells.Item(119, 5).Value = 80
$ws.Cells.Item(119, 6).Value = 20
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 15

